# Append two new daily-charge rows (2025-10-21) to Sheet1, mirroring the
# layout/formatting of the most recent existing rows (100/101).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 102: 四方坪站充电量(kw) -------------------------------------------------
# Copy row 100 (same station) into row 102 so number formats / date format
# (styles s="3" / s="5") are inherited exactly, then overwrite the values.
$ws.Range("A100:Z100").Copy($ws.Range("A102:Z102"))

$ws.Range("A102").Value = 45951
$ws.Range("B102").Value = "四方坪站充电量(kw)"

$row102 = @(857.52099999999984,1466.3820000000003,479.98,539.16899999999998,493.20499999999993,625.346,390.88000000000005,345.15300000000002,185.28000000000003,162.74,172.44899999999998,164.58999999999997,916.23400000000004,1404.0079999999998,490.19,234.20499999999998,224.02,205.00700000000001,68.490000000000009,34.880000000000003,52.459999999999994,99.52,61.91,46.322000000000003)
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "102").Value = $row102[$i]
}

# --- Row 103: 高岭站充电量(kw) --------------------------------------------------
$ws.Range("A101:Z101").Copy($ws.Range("A103:Z103"))

$ws.Range("A103").Value = 45951
$ws.Range("B103").Value = "高岭站充电量(kw)"

$row103 = @(383.12700000000001,717.8839999999999,295.90899999999999,45.890999999999998,19.938000000000002,461.00799999999998,302.74400000000003,160.67800000000003,385.49099999999999,207.84299999999996,110.67400000000001,437.09899999999999,718.36899999999991,743.57199999999978,179.71499999999997,198.86999999999998,130.11500000000001,179.92200000000003,68.594999999999999,62.534000000000006,41.984000000000002,21.210999999999999,34.519000000000005,39.335999999999999)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "103").Value = $row103[$i]
}

# Match the saved selection state from the source edit.
$ws.Range("J106").Select()
